$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.848.26'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.937.14'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.54'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4888'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2952'
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06880'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.28'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '104.78'
$ws.Range('E11').Value = '  -2.93%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07792'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.934.74'
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.348'
$ws.Range('E14').Value = '  -2.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7005'
$ws.Range('E15').Value = '  -1.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '273.05'
$ws.Range('E16').Value = '  -3.63%  '
$ws.Range('D17').Value = '30.882.58'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007726'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.624'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.204.22'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.538'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.842'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.15'
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.61'
$ws.Range('E27').Value = '  -2.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.154'
$ws.Range('E28').Value = '  -3.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1040'
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.391'
$ws.Range('E30').Value = '  -2.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.558'
$ws.Range('E31').Value = '  -1.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.574'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.385'
$ws.Range('E33').Value = '  -2.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04889'
$ws.Range('E34').Value = '  -2.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7602'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.151'
$ws.Range('E36').Value = '  -2.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.735'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02012'
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '79.87'
$ws.Range('E40').Value = '  +6.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.663'
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.504'
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.079'
$ws.Range('E43').Value = '  -4.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9074'
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4443'
$ws.Range('E45').Value = '  -1.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '108.11'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.002'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.751'
$ws.Range('E48').Value = '  -4.64%  '
$ws.Range('D49').Value = '1.000.19'
$ws.Range('E49').Value = '  +1.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1246'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.20'
$ws.Range('E51').Value = '  +0.94%  '
